$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values: force text formatting so that numeric-looking
# strings (e.g. "1.00") are not auto-converted to numbers by Excel, which
# would otherwise drop significant trailing zeros.
$priceUpdates = [ordered]@{
    2 = '26.405.76'
    3 = '1.623.35'
    4 = '1.00'
    5 = '212.47'
    7 = '1.00'
    10 = '18.88'
    12 = '1.850.55'
    13 = '1.646.46'
    17 = '26.408.11'
    19 = '214.39'
    20 = '1.00'
    22 = '6.20'
    24 = '1.96'
    25 = '147.77'
    27 = '0.119'
    29 = '15.53'
    32 = '3.32'
    33 = '2.94'
    36 = '1.212.43'
    38 = '1.00'
    39 = '0.792'
    41 = '2.24'
    44 = '1.759.33'
    45 = '92.61'
    47 = '54.62'
    49 = '0.0510'
    50 = '7.63'
    51 = '0.407'
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

# Column E (Volume(1h)) values: plain text updates (percentages with
# surrounding padding spaces), no special handling needed.
$volumeUpdates = [ordered]@{
    3 = '  -0.80%  '
    4 = '  +0.10%  '
    5 = '  -0.55%  '
    6 = '  +0.81%  '
    7 = '  +0.22%  '
    8 = '  -1.13%  '
    9 = '  -0.06%  '
    10 = '  -0.98%  '
    11 = '  +0.53%  '
    12 = '  -0.77%  '
    13 = '  +0.64%  '
    14 = '  +1.02%  '
    15 = '  -0.83%  '
    16 = '  +1.11%  '
    17 = '  -1.02%  '
    18 = '  +0.60%  '
    19 = '  +2.98%  '
    20 = '  +0.32%  '
    21 = '  -0.97%  '
    22 = '  +1.64%  '
    23 = '  -1.37%  '
    24 = '  +3.48%  '
    25 = '  +1.37%  '
    26 = '  +0.26%  '
    27 = '  -0.77%  '
    28 = '  +2.12%  '
    29 = '  +0.87%  '
    30 = '  -1.97%  '
    31 = '  -1.54%  '
    32 = '  +2.70%  '
    33 = '  -0.43%  '
    34 = '  -1.14%  '
    35 = '  -1.69%  '
    36 = '  +3.81%  '
    37 = '  +3.20%  '
    38 = '  +0.26%  '
    39 = '  -2.67%  '
    40 = '  -0.98%  '
    41 = '  -3.26%  '
    42 = '  -0.37%  '
    43 = '  -0.86%  '
    44 = '  -0.78%  '
    45 = '  +0.20%  '
    46 = '  +1.41%  '
    47 = '  -0.08%  '
    48 = '  -0.30%  '
    49 = '  -0.27%  '
    50 = '  +0.01%  '
    51 = '  -0.58%  '
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}
